# Append 11 new daily rows (2021-04-22 .. 2021-05-02) to the flight-tracking
# table on "Ark1", following the same layout as the existing rows:
#   A = DateTime (text), B = Scheduled flights, C = Tracked flights,
#   D = Percent (=C/B)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 382: 2021-04-22
$ws.Range("A381:D381").Copy()
$ws.Range("A382:D382").PasteSpecial(-4122)
$ws.Cells.Item(382, 1).Value = "2021-04-22"
$ws.Cells.Item(382, 2).Value = 66
$ws.Cells.Item(382, 3).Value = 64
$ws.Cells.Item(382, 4).Formula = "=C382/B382"

# Row 383: 2021-04-23
$ws.Range("A382:D382").Copy()
$ws.Range("A383:D383").PasteSpecial(-4122)
$ws.Cells.Item(383, 1).Value = "2021-04-23"
$ws.Cells.Item(383, 2).Value = 67
$ws.Cells.Item(383, 3).Value = 63
$ws.Cells.Item(383, 4).Formula = "=C383/B383"

# Row 384: 2021-04-24
$ws.Range("A383:D383").Copy()
$ws.Range("A384:D384").PasteSpecial(-4122)
$ws.Cells.Item(384, 1).Value = "2021-04-24"
$ws.Cells.Item(384, 2).Value = 51
$ws.Cells.Item(384, 3).Value = 49
$ws.Cells.Item(384, 4).Formula = "=C384/B384"

# Row 385: 2021-04-25
$ws.Range("A384:D384").Copy()
$ws.Range("A385:D385").PasteSpecial(-4122)
$ws.Cells.Item(385, 1).Value = "2021-04-25"
$ws.Cells.Item(385, 2).Value = 53
$ws.Cells.Item(385, 3).Value = 48
$ws.Cells.Item(385, 4).Formula = "=C385/B385"

# Row 386: 2021-04-26
$ws.Range("A385:D385").Copy()
$ws.Range("A386:D386").PasteSpecial(-4122)
$ws.Cells.Item(386, 1).Value = "2021-04-26"
$ws.Cells.Item(386, 2).Value = 58
$ws.Cells.Item(386, 3).Value = 56
$ws.Cells.Item(386, 4).Formula = "=C386/B386"

# Row 387: 2021-04-27
$ws.Range("A386:D386").Copy()
$ws.Range("A387:D387").PasteSpecial(-4122)
$ws.Cells.Item(387, 1).Value = "2021-04-27"
$ws.Cells.Item(387, 2).Value = 50
$ws.Cells.Item(387, 3).Value = 49
$ws.Cells.Item(387, 4).Formula = "=C387/B387"

# Row 388: 2021-04-28
$ws.Range("A387:D387").Copy()
$ws.Range("A388:D388").PasteSpecial(-4122)
$ws.Cells.Item(388, 1).Value = "2021-04-28"
$ws.Cells.Item(388, 2).Value = 65
$ws.Cells.Item(388, 3).Value = 64
$ws.Cells.Item(388, 4).Formula = "=C388/B388"

# Row 389: 2021-04-29
$ws.Range("A388:D388").Copy()
$ws.Range("A389:D389").PasteSpecial(-4122)
$ws.Cells.Item(389, 1).Value = "2021-04-29"
$ws.Cells.Item(389, 2).Value = 68
$ws.Cells.Item(389, 3).Value = 64
$ws.Cells.Item(389, 4).Formula = "=C389/B389"

# Row 390: 2021-04-30
$ws.Range("A389:D389").Copy()
$ws.Range("A390:D390").PasteSpecial(-4122)
$ws.Cells.Item(390, 1).Value = "2021-04-30"
$ws.Cells.Item(390, 2).Value = 59
$ws.Cells.Item(390, 3).Value = 55
$ws.Cells.Item(390, 4).Formula = "=C390/B390"

# Row 391: 2021-05-01
$ws.Range("A390:D390").Copy()
$ws.Range("A391:D391").PasteSpecial(-4122)
$ws.Cells.Item(391, 1).Value = "2021-05-01"
$ws.Cells.Item(391, 2).Value = 48
$ws.Cells.Item(391, 3).Value = 45
$ws.Cells.Item(391, 4).Formula = "=C391/B391"

# Row 392: 2021-05-02
$ws.Range("A391:D391").Copy()
$ws.Range("A392:D392").PasteSpecial(-4122)
$ws.Cells.Item(392, 1).Value = "2021-05-02"
$ws.Cells.Item(392, 2).Value = 57
$ws.Cells.Item(392, 3).Value = 56
$ws.Cells.Item(392, 4).Formula = "=C392/B392"

$excel.CutCopyMode = $false

# Match the author's final on-screen view: scrolled so row 359 is at the
# top, with the new last row fully selected.
$excel.Goto($ws.Range("A359"), $true) | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 359
$win.ScrollColumn = 1
$ws.Range("A381:XFD381").Select() | Out-Null
